# Insert two new data rows before the current row 185, shifting the existing
# rows 185:300 down to 187:302 (matches the target dimension A1:R302).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("185:186").Insert()

# New row 185: Choclero / Primera record dated 2023-02-03 (serial 44960)
$ws.Range("A185").Value = 5
$ws.Range("B185").Value = "Macroferia Regional de Talca"
$ws.Range("C185").Value = "Maule"
$ws.Range("D185").Value = 44960
$ws.Range("D185").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E185").Value = 7
$ws.Range("F185").Value = 100112024
$ws.Range("G185").Value = "Choclo"
$ws.Range("H185").Value = "Choclero"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 40000
$ws.Range("K185").Value = 300
$ws.Range("L185").Value = 350
$ws.Range("M185").Value = 325
$ws.Range("N185").Value = "`$/unidad"
$ws.Range("O185").Value = "Región del Maule"
$ws.Range("P185").Value = 325
$ws.Range("Q185").Value = 1
$ws.Range("R185").Value = "Hortaliza"

# New row 186: Choclero / Segunda record dated 2023-02-03 (serial 44960)
$ws.Range("A186").Value = 5
$ws.Range("B186").Value = "Macroferia Regional de Talca"
$ws.Range("C186").Value = "Maule"
$ws.Range("D186").Value = 44960
$ws.Range("D186").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E186").Value = 7
$ws.Range("F186").Value = 100112024
$ws.Range("G186").Value = "Choclo"
$ws.Range("H186").Value = "Choclero"
$ws.Range("I186").Value = "Segunda"
$ws.Range("J186").Value = 10000
$ws.Range("K186").Value = 200
$ws.Range("L186").Value = 200
$ws.Range("M186").Value = 200
$ws.Range("N186").Value = "`$/unidad"
$ws.Range("O186").Value = "Región del Maule"
$ws.Range("P186").Value = 200
$ws.Range("Q186").Value = 1
$ws.Range("R186").Value = "Hortaliza"
